# feat: register user and asistence
# Add a new "identityDocument" column/header next to the existing
# user-registration fields (name, lastName, email, password, phoneNumber).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 with the "identityDocument" field name.
$ws.Range("F1").Value = "identityDocument"

# Move/leave the active selection on A6, as in the saved workbook.
$ws.Range("A6").Select() | Out-Null
